$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.027.73"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "1.558.22"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").Value = "'207.35"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("E8").Value = "  +1.33%  "
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").Value = "'0.0591"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("D11").Value = "'0.0863"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "1.780.28"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "1.558.25"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "27.016.13"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "'62.01"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "'216.12"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").Value = "'1.01"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("D23").Value = "'9.23"
$ws.Range("E23").Value = "  +2.61%  "
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("D27").Value = "'14.94"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").Value = "'0.0464"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").Value = "1.401.25"
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("E34").Value = "  +3.24%  "
$ws.Range("E35").Value = "  +3.12%  "
$ws.Range("D36").Value = "'0.961"
$ws.Range("E36").Value = "  +2.82%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").Value = "'0.523"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").Value = "'0.811"
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("D41").Value = "'1.01"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").Value = "'0.991"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("E43").Value = "  +3.32%  "
$ws.Range("D44").Value = "'5.48"
$ws.Range("D45").Value = "'64.01"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").Value = "'1.75"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").Value = "1.694.19"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").Value = "'86.26"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").Value = "'0.0960"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("E51").Value = "  +0.41%  "
